# Update the "想去人数" (interest count) figures in column F for the two
# worksheets that list exhibition rows ("展览" and "全部类型"), matching the
# refreshed counts captured by the scraper run.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3  = 578
    7  = 60
    11 = 4707
    12 = 4492
    13 = 18
    14 = 1
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
